$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "42.383.13"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").Value = "2.283.67"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws.Range("D5") "157.53"
$ws.Range("E5").Value = "  +15,645.30%  "
Set-TextValue $ws.Range("D6") "306.28"
$ws.Range("E6").Value = "  +1.09%  "
Set-TextValue $ws.Range("D7") "95.58"
$ws.Range("E7").Value = "  +4.94%  "
Set-TextValue $ws.Range("D8") "0.530"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("E9").Value = "  +0.01%  "
Set-TextValue $ws.Range("D10") "0.493"
$ws.Range("E10").Value = "  +2.87%  "
Set-TextValue $ws.Range("D11") "35.72"
$ws.Range("E11").Value = "  +11.46%  "
Set-TextValue $ws.Range("D12") "0.0802"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("E13").Value = "  -2.01%  "
Set-TextValue $ws.Range("D14") "6.70"
$ws.Range("E14").Value = "  +2.10%  "
$ws.Range("D15").Value = "2.637.51"
$ws.Range("E15").Value = "  +1.14%  "
Set-TextValue $ws.Range("D16") "14.49"
$ws.Range("E16").Value = "  +2.67%  "
$ws.Range("D17").Value = "2.294.18"
$ws.Range("E17").Value = "  +1.93%  "
Set-TextValue $ws.Range("D18") "0.799"
$ws.Range("E18").Value = "  +5.50%  "
$ws.Range("D19").Value = "42.293.54"
$ws.Range("E19").Value = "  +1.76%  "
Set-TextValue $ws.Range("D20") "12.68"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").Value = "0.0₃0917"
$ws.Range("E21").Value = "  +1.70%  "
Set-TextValue $ws.Range("D22") "5.99"
$ws.Range("E22").Value = "  +1.79%  "
Set-TextValue $ws.Range("D23") "67.87"
$ws.Range("E23").Value = "  +2.05%  "
Set-TextValue $ws.Range("D24") "242.54"
$ws.Range("E24").Value = "  +0.96%  "
Set-TextValue $ws.Range("D25") "2.60"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("E26").Value = "  +1.79%  "
$ws.Range("E27").Value = "  -0.23%  "
Set-TextValue $ws.Range("D28") "23.91"
$ws.Range("E28").Value = "  -0.39%  "
Set-TextValue $ws.Range("D29") "35.75"
$ws.Range("E29").Value = "  +4.12%  "
Set-TextValue $ws.Range("D30") "9.56"
$ws.Range("E30").Value = "  +0.80%  "
Set-TextValue $ws.Range("D31") "2.10"
$ws.Range("E31").Value = "  +1.89%  "
Set-TextValue $ws.Range("D32") "160.93"
$ws.Range("E32").Value = "  -0.15%  "
Set-TextValue $ws.Range("D33") "5.30"
$ws.Range("E33").Value = "  +3.31%  "
$ws.Range("E34").Value = "  +0.10%  "
Set-TextValue $ws.Range("D35") "0.0752"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("E36").Value = "  +3.20%  "
$ws.Range("E37").Value = "  +4.99%  "
Set-TextValue $ws.Range("D38") "17.16"
$ws.Range("E38").Value = "  +3.69%  "
$ws.Range("E40").Value = "  +3.38%  "
Set-TextValue $ws.Range("D42") "4.14"
$ws.Range("E42").Value = "  +6.01%  "
$ws.Range("D43").Value = "2.004.60"
Set-TextValue $ws.Range("D44") "2.30"
$ws.Range("E44").Value = "  +12.62%  "
Set-TextValue $ws.Range("D45") "19.20"
$ws.Range("E45").Value = "  -1.39%  "
Set-TextValue $ws.Range("D46") "0.0284"
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D47") "10.16"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D48") "2.99"
$ws.Range("E48").Value = "  +5.22%  "
Set-TextValue $ws.Range("D49") "53.48"
$ws.Range("E49").Value = "  +4.09%  "
$ws.Range("E50").Value = "  +1.80%  "
Set-TextValue $ws.Range("D51") "72.66"
$ws.Range("E51").Value = "  +0.65%  "
